# Apply the "Added mods for NGR filter selects and 50 Ohm terminations" edit.
#
# Summary of content changes:
#  - SerDes sheet, column L (rows 20-49): swap the "Next*"/"Prev*" labels
#    between the two filter-select groups (L20-L32 <-> L36-L49).
#  - Selection/active-sheet state changes: SerDes becomes the active sheet
#    (tabSelected) with the current selection on L36:L49, and Sheet1 is no
#    longer the tab-selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SerDes")

# --- Update column L labels (swap Next*/Prev* between the two groups) ---
$ws.Range("L20").Value = "PrevClk"
$ws.Range("L23").Value = "PrevD1"
$ws.Range("L25").Value = "PrevD2"
$ws.Range("L27").Value = "PrevD3"
$ws.Range("L29").Value = "PrevD4"
$ws.Range("L32").Value = "PrevD0"

$ws.Range("L36").Value = "NextClk"
$ws.Range("L39").Value = "NextD4"
$ws.Range("L41").Value = "NextD3"
$ws.Range("L44").Value = "NextD2"
$ws.Range("L46").Value = "NextD1"
$ws.Range("L49").Value = "NextD0"

# --- Update view/selection state ---
# Make SerDes the active sheet (tabSelected), scroll it so row 16 is at the
# top, and select L36:L49 as the current selection on that sheet (matches
# active cell L36).
$ws.Activate()

$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$win.Left = 26745
$win.Top = 13695
$win.Width = 24135
$win.Height = 15600

$ws.Range("L36:L49").Select()
